$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (columns A-H): anchor words under the "negative" header ---
# The "fraud" row is gone (larger re-run dataset has only 4 negative anchors now),
# so fully clear what used to be the 5th data row (row 7) including its style.
$ws.Range("A7:H7").Clear()

$ws.Range("A3").Value = 'crude'
$ws.Range("B3").Value = 0.9117647058823529
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 31
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 3

$ws.Range("A4").Value = 'crisis'
$ws.Range("B4").Value = 0.6404109589041096
$ws.Range("C4").Value = 187
$ws.Range("D4").Value = 187
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 105

$ws.Range("A5").Value = 'panic'
$ws.Range("B5").Value = 0.1937984496124031
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 416

$ws.Range("A6").Value = 'sc'
$ws.Range("B6").Value = 0.1534391534391534
$ws.Range("C6").Value = 29
$ws.Range("D6").Value = 29
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 160

# --- Table 2 (columns J-Q): anchor words under the "positive" header ---
# Still 25 data rows (3-27), but word list / stats changed with the larger re-run.

$ws.Range("J3").Value = 'love'
$ws.Range("K3").Value = 0.9782608695652174
$ws.Range("L3").Value = 45
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 1

$ws.Range("J4").Value = 'interesting'
$ws.Range("K4").Value = 0.9696969696969697
$ws.Range("L4").Value = 32
$ws.Range("M4").Value = 32
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 1

$ws.Range("J5").Value = 'best'
$ws.Range("K5").Value = 0.9152542372881356
$ws.Range("L5").Value = 54
$ws.Range("M5").Value = 54
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 5

$ws.Range("J6").Value = 'great'
$ws.Range("K6").Value = 0.875
$ws.Range("L6").Value = 98
$ws.Range("M6").Value = 98
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 14

$ws.Range("J7").Value = 'thanks'
$ws.Range("K7").Value = 0.8414634146341463
$ws.Range("L7").Value = 69
$ws.Range("M7").Value = 69
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 13

$ws.Range("J8").Value = 'special'
$ws.Range("K8").Value = 0.8333333333333334
$ws.Range("L8").Value = 30
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 6

$ws.Range("J9").Value = 'positive'
$ws.Range("K9").Value = 0.8275862068965517
$ws.Range("L9").Value = 48
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 10

$ws.Range("J10").Value = 'thank'
$ws.Range("K10").Value = 0.8203125
$ws.Range("L10").Value = 105
$ws.Range("M10").Value = 105
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 23

$ws.Range("J11").Value = 'free'
$ws.Range("K11").Value = 0.7916666666666666
$ws.Range("L11").Value = 95
$ws.Range("M11").Value = 95
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 25

$ws.Range("J12").Value = 'safe'
$ws.Range("K12").Value = 0.7323943661971831
$ws.Range("L12").Value = 104
$ws.Range("M12").Value = 104
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 38

$ws.Range("J13").Value = 'safety'
$ws.Range("K13").Value = 0.7254901960784313
$ws.Range("L13").Value = 37
$ws.Range("M13").Value = 37
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 14

$ws.Range("J14").Value = 'support'
$ws.Range("K14").Value = 0.7075471698113207
$ws.Range("L14").Value = 75
$ws.Range("M14").Value = 75
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 31

$ws.Range("J15").Value = 'good'
$ws.Range("K15").Value = 0.7
$ws.Range("L15").Value = 112
$ws.Range("M15").Value = 112
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 48

$ws.Range("J16").Value = 'confidence'
$ws.Range("K16").Value = 0.6944444444444444
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 11

$ws.Range("J17").Value = 'heroes'
$ws.Range("K17").Value = 0.6808510638297872
$ws.Range("L17").Value = 32
$ws.Range("M17").Value = 32
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 15

$ws.Range("J18").Value = 'well'
$ws.Range("K18").Value = 0.648936170212766
$ws.Range("L18").Value = 61
$ws.Range("M18").Value = 61
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 33

$ws.Range("J19").Value = 'better'
$ws.Range("K19").Value = 0.6349206349206349
$ws.Range("L19").Value = 40
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 23

$ws.Range("J20").Value = 'relief'
$ws.Range("K20").Value = 0.56
$ws.Range("L20").Value = 28
$ws.Range("M20").Value = 28
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 22

$ws.Range("J21").Value = 'fresh'
$ws.Range("K21").Value = 0.5208333333333334
$ws.Range("L21").Value = 25
$ws.Range("M21").Value = 25
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 23

$ws.Range("J22").Value = 'hand'
$ws.Range("K22").Value = 0.5143603133159269
$ws.Range("L22").Value = 197
$ws.Range("M22").Value = 197
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 186

$ws.Range("J23").Value = 'care'
$ws.Range("K23").Value = 0.4719101123595505
$ws.Range("L23").Value = 42
$ws.Range("M23").Value = 42
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 47

$ws.Range("J24").Value = 'like'
$ws.Range("K24").Value = 0.4294117647058823
$ws.Range("L24").Value = 146
$ws.Range("M24").Value = 146
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 194

$ws.Range("J25").Value = 'help'
$ws.Range("K25").Value = 0.4101694915254237
$ws.Range("L25").Value = 121
$ws.Range("M25").Value = 121
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 174

$ws.Range("J26").Value = 'increase'
$ws.Range("K26").Value = 0.3717948717948718
$ws.Range("L26").Value = 29
$ws.Range("M26").Value = 29
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 49

$ws.Range("J27").Value = 'please'
$ws.Range("K27").Value = 0.301255230125523
$ws.Range("L27").Value = 72
$ws.Range("M27").Value = 72
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 167
